# feat: expand arcane depths ui integration
#
# Rebuilds rows 7..15 of the rooms sheet: the original 3 rows (7-9) are
# replaced with 9 new rows covering the extra Elite/Boss/Event/Combat
# entries added for the arcane-depths chapters, and the sheet's used
# range grows from A1:N9 to A1:N15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-like columns that must stay text instead of being auto-coerced to
# numbers by Excel (leading zeros in A-D, purely-numeric relic ids in M).
$ws.Range("A7:D15").NumberFormat = "@"
$ws.Range("M7:M15").NumberFormat = "@"

$rows = @(
    @{ Row=7;  A="50"; B="04"; C="01"; D="00"; E=1; F="Elite";   G="熔核军械库";     H=10010100; I="对抗熔核巨灵与其护卫。";             J=40010200; K=$null;     L="Relic";    M="60010000";     N="精英常驻熔火护盾，每回合爆发灼烧波" }
    @{ Row=8;  A="50"; B="06"; C="02"; D="00"; E=2; F="Boss";    G="灰烬之主燃殿";   H=10010100; I="灰烬之主统御火种，终焉决战即将开启。"; J=40010300; K=$null;     L="Relic";    M="60010000";     N="Boss首次登场时施加全场灼烧" }
    @{ Row=9;  A="50"; B="01"; C="02"; D="00"; E=2; F="Combat";  G="蔓生之巢";       H=10020100; I="藤蔓潮汐的核心据点。";               J=40020300; K=$null;     L="Resource"; M="Provision:18"; N="初始时所有敌人持有藤蔓护盾" }
    @{ Row=10; A="50"; B="02"; C="02"; D="00"; E=2; F="Event";   G="深根共鸣";       H=10020100; I="旅者请求帮助，分享绿色共鸣。";       J=$null;    K=51020000;  L="Resource"; M="Provision:15"; N="成功则提升当前层的治愈效率" }
    @{ Row=11; A="50"; B="05"; C="01"; D="00"; E=1; F="Rest";    G="翠息静室";       H=10020100; I="通过深绿脉络恢复与净化。";           J=$null;    K=$null;     L="Trait";    M="Synergy";       N="移除任意一个减益" }
    @{ Row=12; A="50"; B="06"; C="03"; D="00"; E=3; F="Boss";    G="藤界之心";       H=10020100; I="削弱深根巨树的束缚并阻止其觉醒。";   J=40020300; K=$null;     L="Relic";    M="60020000";     N="Boss召唤孢子爪牙并周期性缠绕全队" }
    @{ Row=13; A="50"; B="01"; C="03"; D="00"; E=3; F="Combat";  G="星火熔层";       H=10030100; I="熔炉中游走的星火构装体。";           J=40030100; K=$null;     L="Resource"; M="Arcane:16";    N="地图施加减速并周期性落下星火" }
    @{ Row=14; A="50"; B="04"; C="02"; D="00"; E=2; F="Elite";   G="星界锻卫";       H=10030100; I="星界傀儡组成的防线。";               J=40030100; K=$null;     L="Relic";    M="60030000";     N="精英拥有反射护盾与星能回复" }
    @{ Row=15; A="50"; B="06"; C="01"; D="00"; E=1; F="Boss";    G="熔炉主宰之厅";   H=10030100; I="对决星界熔炉的主宰。";               J=40030300; K=$null;     L="Relic";    M="60030000";     N="Boss阶段转换召唤星火残迹" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I

    if ($r.J -ne $null) {
        $ws.Cells.Item($n, 10).Value = $r.J
    } else {
        $ws.Cells.Item($n, 10).ClearContents()
    }

    if ($r.K -ne $null) {
        $ws.Cells.Item($n, 11).Value = $r.K
    } else {
        $ws.Cells.Item($n, 11).ClearContents()
    }

    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
    $ws.Cells.Item($n, 14).Value = $r.N
}
